# Add new rodent-trap survey rows (CC location, 5/15-5/16/2020) to the
# "Records" sheet, matching the existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Records")

# New data rows to append (row, Location, DateSerial, Trap#, Replaced#, Sex)
$rows = @(
    @{ Row = 26; A = "CC"; B = 43966; D = 3091; E = 3042; F = $null },
    @{ Row = 27; A = "CC"; B = 43966; D = 3054; E = 3018; F = "F" },
    @{ Row = 28; A = "CC"; B = 43966; D = 3046; E = 3041; F = "F" },
    @{ Row = 29; A = "CC"; B = 43966; D = 3085; E = 3084; F = "J" },
    @{ Row = 30; A = "CC"; B = 43966; D = 3078; E = 3087; F = "F" },
    @{ Row = 31; A = "CC"; B = 43966; D = 3012; E = 3039; F = "F" },
    @{ Row = 32; A = "CC"; B = 43966; D = 3098; E = 3044; F = "?" },
    @{ Row = 33; A = "CC"; B = 43967; D = 3087; E = 3003; F = "M" },
    @{ Row = 34; A = "CC"; B = 43967; D = 3084; E = 3038; F = "J" },
    @{ Row = 35; A = "CC"; B = 43967; D = 3041; E = 3026; F = "F?" },
    @{ Row = 36; A = $null; B = $null; D = 3052; E = 3037; F = "J" }
)

# Use the existing dated cell's format (short date, style index 1) as the
# template for the new Date column entries.
$ws.Range("B25").Copy()

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($r.A) {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
    }

    if ($r.B) {
        $ws.Cells.Item($rowNum, 2).PasteSpecial(-4122)
        $ws.Cells.Item($rowNum, 2).Value = $r.B
    }

    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E

    if ($r.F) {
        $ws.Cells.Item($rowNum, 6).Value = $r.F
    }
}

$excel.CutCopyMode = 0

# Update the view to scroll to / select the newly entered rows.
$ws.Activate()
$excel.Goto($ws.Range("A27"), $true)
$ws.Range("A28:A35").Select()
